$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "Khoni"
$ws.Name = "Khoni"

# Rows 6 ("Urban") and 7 ("Rural") had their yearly figures suppressed:
#   - the 2010 column (B) now shows "..." (a literal three-dot placeholder)
#   - all other year columns (C:O, 2011-2023) now show "…" (the existing
#     ellipsis already used for confidential/unavailable data)
# The "Total" row (5) keeps its real figures untouched.
$ws.Range("B6").Value = "..."
$ws.Range("C6:O6").Value = "…"

$ws.Range("B7").Value = "..."
$ws.Range("C7:O7").Value = "…"

# Row 8 was an entirely blank spacer row between the data table and the
# footnote; it was removed, shifting the footnote up from row 9 to row 8.
$ws.Rows.Item(8).Delete()
